$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 17, shifting rows 17:100 down to 18:101
$ws.Rows("17:17").Insert()

# Populate the new row 17 with its data (matches the other rows' constant
# columns A, B, C, E, F, G, I, Q, R, plus the new record's own values)
$ws.Range("A17").Value = 4
$ws.Range("B17").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C17").Value = 'Los Lagos'
$ws.Range("D17").Value = '12/17/2021'
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 100112022
$ws.Range("G17").Value = 'Arveja Verde'
$ws.Range("H17").Value = 'Sin especificar'
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 90
$ws.Range("K17").Value = 20000
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = 20000
$ws.Range("N17").Value = '$/saco 25 kilos'
$ws.Range("O17").Value = 'Región de La Araucanía'
$ws.Range("P17").Value = 800
$ws.Range("Q17").Value = 25
$ws.Range("R17").Value = 'Hortaliza'
